# Update release version strings on the 'About' sheet
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$about.Range("A2").Value = 'Version: mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)'
$about.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for V. I. Lenin (Kazakhstan) Coal Mine, Kazakhstan, M1438, version ''mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# Update the build_version column (S) on the data sheet for existing rows 2-13
$data = $wb.Worksheets.Item("Boundaries and methane sources")
$newVersion = 'mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)'
for ($r = 2; $r -le 13; $r++) {
    $data.Cells.Item($r, 19).Value = $newVersion
}

# Append two new point-feature rows (14 and 15)
$data.Cells.Item(14, 1).Value = 'V. I. Lenin (Kazakhstan) Coal Mine, Kazakhstan, M1438'
$data.Cells.Item(14, 2).Value = 'M1438.P14'
$data.Cells.Item(14, 3).Value = 'ventilation system'
$data.Cells.Item(14, 4).Value = 'vent'
$data.Cells.Item(14, 5).Value = 45078
$data.Cells.Item(14, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$data.Cells.Item(14, 6).Value = 'Identified visually from Google Earth Pro satellite imagery.'
$data.Cells.Item(14, 7).Value = 'ventilation shaft'
$data.Cells.Item(14, 8).Value = 'extracted'
$data.Cells.Item(14, 9).Value = 'M1438'
$data.Cells.Item(14, 10).Value = 'Qarmet JSC [100%]'
$data.Cells.Item(14, 12).Value = 'Qazaqstan Steel Group LLP'
$data.Cells.Item(14, 13).Value = 'https://www.gem.wiki/V._I._Lenin_(Kazakhstan)_coal_mine'
$data.Cells.Item(14, 15).Value = 'Met'
$data.Cells.Item(14, 16).Value = 'V. I. Lenin (Kazakhstan) Coal Mine'
$data.Cells.Item(14, 17).Value = 'Kazakhstan'
$data.Cells.Item(14, 18).Value = 'Apr 15, 2025'
$data.Cells.Item(14, 19).Value = 'mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)'
$data.Cells.Item(14, 20).Value = 'POINT (72.530966 49.746095)'

$data.Cells.Item(15, 1).Value = 'V. I. Lenin (Kazakhstan) Coal Mine, Kazakhstan, M1438'
$data.Cells.Item(15, 2).Value = 'M1438.P15'
$data.Cells.Item(15, 3).Value = 'degasification system'
$data.Cells.Item(15, 4).Value = 'drainage station'
$data.Cells.Item(15, 5).Value = 45078
$data.Cells.Item(15, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$data.Cells.Item(15, 6).Value = 'Identified visually from Google Earth Pro satellite imagery.'
$data.Cells.Item(15, 7).Value = 'possible gas drainage station'
$data.Cells.Item(15, 8).Value = 'extracted'
$data.Cells.Item(15, 9).Value = 'M1438'
$data.Cells.Item(15, 10).Value = 'Qarmet JSC [100%]'
$data.Cells.Item(15, 12).Value = 'Qazaqstan Steel Group LLP'
$data.Cells.Item(15, 13).Value = 'https://www.gem.wiki/V._I._Lenin_(Kazakhstan)_coal_mine'
$data.Cells.Item(15, 15).Value = 'Met'
$data.Cells.Item(15, 16).Value = 'V. I. Lenin (Kazakhstan) Coal Mine'
$data.Cells.Item(15, 17).Value = 'Kazakhstan'
$data.Cells.Item(15, 18).Value = 'Apr 15, 2025'
$data.Cells.Item(15, 19).Value = 'mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)'
$data.Cells.Item(15, 20).Value = 'POINT (72.513121 49.735834)'

Write-Host "Version strings and new rows updated."
